$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a new column at C, shifting existing columns C.. to the right
$ws.Columns.Item(3).Insert()

# Populate the new "rows" header column
$ws.Cells.Item(1, 3).Value = "rows"
$ws.Cells.Item(4, 3).Value = 12

$ws.Range("D5").Select()
